# Apply "agregado de password en import xlsx" changes
$wb = $excel.ActiveWorkbook

# --- Sheet "Instrucciones": update client field docs and notes ---
$wsInstr = $wb.Worksheets.Item("Instrucciones")

# Row 7: address is now optional instead of required
$wsInstr.Cells.Item(7, 1).Value = "- address: Dirección completa (opcional)"

# Row 8: replace "notes" field doc with new "password" field doc
$wsInstr.Cells.Item(8, 1).Value = "- password: Contraseña del cliente (obligatorio)"

# Rows 21-27: new/shifted notes at the bottom of the sheet
$wsInstr.Cells.Item(21, 1).Value = "- Campos obligatorios para clientes: name, email, phone, password"
$wsInstr.Cells.Item(22, 1).Value = "- Campo opcional para clientes: address (puede dejarse vacío)"
$wsInstr.Cells.Item(23, 1).Value = "- Los emails deben ser únicos y válidos"
$wsInstr.Cells.Item(24, 1).Value = "- Las placas deben ser únicas"
$wsInstr.Cells.Item(25, 1).Value = "- El año debe ser un número válido"
$wsInstr.Cells.Item(26, 1).Value = "- Las contraseñas se almacenarán tal como se escriban"
$wsInstr.Cells.Item(27, 1).Value = "- Asegúrese de que el email del cliente existe antes de asignar vehículos"

# --- Sheet "Clientes": replace "notes" column with "password" column ---
$wsClientes = $wb.Worksheets.Item("Clientes")

$wsClientes.Cells.Item(1, 5).Value = "password"
$wsClientes.Cells.Item(2, 5).Value = "mipassword123"
$wsClientes.Cells.Item(3, 5).Value = "password456"

# Column E width shrinks from 30 to 15
# (Excel stores width in a pixel-rounded unit; 14.15 round-trips to exactly 15)
$wsClientes.Columns.Item(5).ColumnWidth = 14.15
